# Auto-generated Excel COM-interop script
# Applies numeric corrections to Leve profit calculation columns (H,I,J,K,L,M,N)
# across multiple sheets, as described in the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 978614.4
$ws.Range("J17").Value = 1097921
$ws.Range("L17").Value = 3293763
$ws.Range("N17").Value = -3294099
# Row 28
$ws.Range("H28").Value = 752.73334
$ws.Range("I28").Value = 476.9091
$ws.Range("K28").Value = 476.9091
$ws.Range("M28").Value = 8.090899999999976
# Row 38
$ws.Range("H38").Value = 3341.4814
$ws.Range("I38").Value = 2310.3572
$ws.Range("J38").Value = 4451.923
$ws.Range("K38").Value = 6931.071599999999
$ws.Range("L38").Value = 13355.769
$ws.Range("M38").Value = -6559.071599999999
$ws.Range("N38").Value = -14099.769
# Row 39
$ws.Range("H39").Value = 706.94446
$ws.Range("I39").Value = 48.333332
$ws.Range("K39").Value = 144.999996
$ws.Range("M39").Value = 151.000004
# Row 40
$ws.Range("H40").Value = 2600.3333
$ws.Range("I40").Value = 2600.3333
$ws.Range("K40").Value = 2600.3333
$ws.Range("M40").Value = -2425.3333
# Row 51
$ws.Range("H51").Value = 4970.619
$ws.Range("J51").Value = 5049.15
$ws.Range("L51").Value = 5049.15
$ws.Range("N51").Value = -6017.15
# Row 61
$ws.Range("H61").Value = 947
$ws.Range("I61").Value = 947
$ws.Range("K61").Value = 2841
$ws.Range("M61").Value = -2669
# Row 96
$ws.Range("H96").Value = 502.42856
$ws.Range("I96").Value = 502.42856
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1507.28568
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -134.28568
$ws.Range("N96").ClearContents()
# Row 111
$ws.Range("H111").Value = 2946.818
$ws.Range("I111").Value = 5094
$ws.Range("J111").Value = 1157.5
$ws.Range("K111").Value = 15282
$ws.Range("L111").Value = 3472.5
$ws.Range("M111").Value = -12215
$ws.Range("N111").Value = -9606.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3101.6365
$ws.Range("I2").Value = 1617.4286
$ws.Range("J2").Value = 5699
$ws.Range("K2").Value = 1617.4286
$ws.Range("L2").Value = 5699
$ws.Range("M2").Value = -1504.4286
$ws.Range("N2").Value = -5925
# Row 32
$ws.Range("H32").Value = 4974.7954
$ws.Range("I32").Value = 4509.875
$ws.Range("J32").Value = 9624
$ws.Range("K32").Value = 4509.875
$ws.Range("L32").Value = 9624
$ws.Range("M32").Value = -4222.875
$ws.Range("N32").Value = -10198
# Row 45
$ws.Range("H45").Value = 3791.8
$ws.Range("I45").Value = 3317.5
$ws.Range("K45").Value = 3317.5
$ws.Range("M45").Value = -2940.5
# Row 74
$ws.Range("H74").Value = 45459292
$ws.Range("I74").Value = 50004696
$ws.Range("J74").Value = 5250
$ws.Range("K74").Value = 50004696
$ws.Range("L74").Value = 5250
$ws.Range("M74").Value = -50003822
$ws.Range("N74").Value = -6998
# Row 77
$ws.Range("H77").Value = 45459292
$ws.Range("I77").Value = 50004696
$ws.Range("J77").Value = 5250
$ws.Range("K77").Value = 250023480
$ws.Range("L77").Value = 26250
$ws.Range("M77").Value = -250019112
$ws.Range("N77").Value = -34986
# Row 116
$ws.Range("H116").Value = 3101.6365
$ws.Range("I116").Value = 1617.4286
$ws.Range("J116").Value = 5699
$ws.Range("K116").Value = 1617.4286
$ws.Range("L116").Value = 5699
$ws.Range("M116").Value = 676.5714
$ws.Range("N116").Value = -10287
# Row 132
$ws.Range("H132").Value = 37101348
$ws.Range("I132").Value = 15258.333
$ws.Range("J132").Value = 166902670
$ws.Range("K132").Value = 45774.999
$ws.Range("L132").Value = 500708010
$ws.Range("M132").Value = -43244.999
$ws.Range("N132").Value = -500713070

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3101.6365
$ws.Range("I3").Value = 1617.4286
$ws.Range("J3").Value = 5699
$ws.Range("K3").Value = 1617.4286
$ws.Range("L3").Value = 5699
$ws.Range("M3").Value = -1503.4286
$ws.Range("N3").Value = -5927
# Row 86
$ws.Range("H86").Value = 62187.875
$ws.Range("I86").Value = 48001
$ws.Range("J86").Value = 104748.5
$ws.Range("K86").Value = 48001
$ws.Range("L86").Value = 104748.5
$ws.Range("M86").Value = -46878
$ws.Range("N86").Value = -106994.5
# Row 89
$ws.Range("H89").Value = 62187.875
$ws.Range("I89").Value = 48001
$ws.Range("J89").Value = 104748.5
$ws.Range("K89").Value = 240005
$ws.Range("L89").Value = 523742.5
$ws.Range("M89").Value = -234389
$ws.Range("N89").Value = -534974.5
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 103
$ws.Range("H103").Value = 65700
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 65700
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 65700
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -68044
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 107
$ws.Range("H107").Value = 3522.2
$ws.Range("I107").Value = 3691.3333
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3691.3333
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -1771.3333
$ws.Range("N107").Value = -5840
# Row 108
$ws.Range("H108").Value = 110500
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 134
$ws.Range("H134").Value = 12000
$ws.Range("I134").Value = 12000
$ws.Range("K134").Value = 36000
$ws.Range("M134").Value = -33465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1450.4615
$ws.Range("I22").Value = 339.55554
$ws.Range("K22").Value = 339.55554
$ws.Range("M22").Value = 10.44445999999999
# Row 109
$ws.Range("H109").Value = 46833
$ws.Range("I109").Value = 54999.5
$ws.Range("J109").Value = 42749.75
$ws.Range("K109").Value = 54999.5
$ws.Range("L109").Value = 42749.75
$ws.Range("M109").Value = -53959.5
$ws.Range("N109").Value = -44829.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 342.86957
$ws.Range("J12").Value = 405.1875
$ws.Range("L12").Value = 1215.5625
$ws.Range("N12").Value = -1561.5625
# Row 98
$ws.Range("H98").Value = 646.3333
$ws.Range("J98").Value = 652.125
$ws.Range("L98").Value = 1956.375
$ws.Range("N98").Value = -4952.375
# Row 100
$ws.Range("H100").Value = 2750
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 9000
$ws.Range("N100").Value = -10622

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 27333.334
# Row 132
$ws.Range("H132").Value = 3270.6667
$ws.Range("I132").Value = 3270.6667
$ws.Range("K132").Value = 9812.000100000001
$ws.Range("M132").Value = -7282.000100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5726.294
$ws.Range("I7").Value = 4493.875
$ws.Range("K7").Value = 4493.875
$ws.Range("M7").Value = -4381.875
# Row 46
$ws.Range("H46").Value = 1446.2572
$ws.Range("I46").Value = 986.931
$ws.Range("J46").Value = 3666.3333
$ws.Range("K46").Value = 986.931
$ws.Range("L46").Value = 3666.3333
$ws.Range("M46").Value = -798.931
$ws.Range("N46").Value = -4042.3333
# Row 68
$ws.Range("H68").Value = 5089.364
$ws.Range("I68").Value = 3998.6667
$ws.Range("J68").Value = 9997.5
$ws.Range("K68").Value = 3998.6667
$ws.Range("L68").Value = 9997.5
$ws.Range("M68").Value = -3249.6667
$ws.Range("N68").Value = -11495.5
# Row 71
$ws.Range("H71").Value = 5089.364
$ws.Range("I71").Value = 3998.6667
$ws.Range("J71").Value = 9997.5
$ws.Range("K71").Value = 19993.3335
$ws.Range("L71").Value = 49987.5
$ws.Range("M71").Value = -16249.3335
$ws.Range("N71").Value = -57475.5
# Row 82
$ws.Range("H82").Value = 2162.4
$ws.Range("I82").Value = 2157.5
$ws.Range("J82").Value = 2165.6667
$ws.Range("K82").Value = 2157.5
$ws.Range("L82").Value = 2165.6667
$ws.Range("M82").Value = -1796.5
$ws.Range("N82").Value = -2887.6667
# Row 85
$ws.Range("H85").Value = 2162.4
$ws.Range("I85").Value = 2157.5
$ws.Range("J85").Value = 2165.6667
$ws.Range("K85").Value = 2157.5
$ws.Range("L85").Value = 2165.6667
$ws.Range("M85").Value = -909.5
$ws.Range("N85").Value = -4661.6667
# Row 116
$ws.Range("H116").Value = 206500
$ws.Range("J116").Value = 206500
$ws.Range("L116").Value = 206500
$ws.Range("N116").Value = -215678
# Row 120
$ws.Range("H120").Value = 56698
$ws.Range("J120").Value = 56698
$ws.Range("L120").Value = 56698
$ws.Range("N120").Value = -66374
# Row 126
$ws.Range("H126").Value = 5726.294
$ws.Range("I126").Value = 4493.875
$ws.Range("K126").Value = 13481.625
$ws.Range("M126").Value = -11011.625
# Row 132
$ws.Range("H132").Value = 3464.9714
$ws.Range("I132").Value = 3321.077
$ws.Range("J132").Value = 3880.6667
$ws.Range("K132").Value = 9963.231
$ws.Range("L132").Value = 11642.0001
$ws.Range("M132").Value = -7433.231
$ws.Range("N132").Value = -16702.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 5050.3
$ws.Range("I126").Value = 6161.8
$ws.Range("J126").Value = 1715.8
$ws.Range("K126").Value = 18485.4
$ws.Range("L126").Value = 5147.4
$ws.Range("M126").Value = -16015.4
$ws.Range("N126").Value = -10087.4
